$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ItemDatas")

# --- Phase 1: seed new shared strings in the exact order the original file introduced them ---
# (so the underlying sharedStrings table comes out index-for-index identical)
$ws.Cells.Item(75, "B").Value = '딸기'
$ws.Cells.Item(76, "B").Value = '딸기 씨'
$ws.Cells.Item(77, "B").Value = '대파'
$ws.Cells.Item(78, "B").Value = '대파 씨'
$ws.Cells.Item(79, "B").Value = '감자'
$ws.Cells.Item(80, "B").Value = '감자 씨'
$ws.Cells.Item(81, "B").Value = '양파'
$ws.Cells.Item(82, "B").Value = '양파 씨'
$ws.Cells.Item(83, "B").Value = '당근'
$ws.Cells.Item(84, "B").Value = '당근 씨'
$ws.Cells.Item(85, "B").Value = '블루베리'
$ws.Cells.Item(86, "B").Value = '블루베리 씨'
$ws.Cells.Item(87, "B").Value = '무'
$ws.Cells.Item(88, "B").Value = '무 씨'
$ws.Cells.Item(89, "B").Value = '양배추'
$ws.Cells.Item(90, "B").Value = '양배추 씨'
$ws.Cells.Item(91, "B").Value = '콜리플라워'
$ws.Cells.Item(92, "B").Value = '콜리플라워 씨'
$ws.Cells.Item(93, "B").Value = '밀'
$ws.Cells.Item(94, "B").Value = '밀 씨'
$ws.Cells.Item(95, "B").Value = '브로콜리'
$ws.Cells.Item(96, "B").Value = '브로콜리 씨'
$ws.Cells.Item(75, "L").Value = 'Strawberry'
$ws.Cells.Item(76, "L").Value = 'Strawberry_Seed'
$ws.Cells.Item(77, "L").Value = 'Greenonion'
$ws.Cells.Item(81, "L").Value = 'Onion'
$ws.Cells.Item(83, "L").Value = 'Carrot'
$ws.Cells.Item(85, "L").Value = 'Blueberry'
$ws.Cells.Item(87, "L").Value = 'Radish'
$ws.Cells.Item(89, "L").Value = 'Cabbage'
$ws.Cells.Item(91, "L").Value = 'Cauliflower'
$ws.Cells.Item(93, "L").Value = 'Wheat'
$ws.Cells.Item(95, "L").Value = 'Broccoli'
$ws.Cells.Item(79, "L").Value = 'Potato'
$ws.Cells.Item(78, "L").Value = 'Greenonion_Seed'
$ws.Cells.Item(80, "L").Value = 'Potato_Seed'
$ws.Cells.Item(82, "L").Value = 'Onion_Seed'
$ws.Cells.Item(84, "L").Value = 'Carrot_Seed'
$ws.Cells.Item(86, "L").Value = 'Blueberry_Seed'
$ws.Cells.Item(88, "L").Value = 'Radish_Seed'
$ws.Cells.Item(90, "L").Value = 'Cabbage_Seed'
$ws.Cells.Item(94, "L").Value = 'Wheat_Seed'
$ws.Cells.Item(92, "L").Value = 'Cauliflower_Seed'
$ws.Cells.Item(96, "L").Value = 'Broccoli_Seed'
$ws.Cells.Item(97, "B").Value = '토마토'
$ws.Cells.Item(98, "B").Value = '토마토 씨'
$ws.Cells.Item(99, "B").Value = '해바라기'
$ws.Cells.Item(100, "B").Value = '해바라기 씨'
$ws.Cells.Item(101, "B").Value = '고추'
$ws.Cells.Item(102, "B").Value = '고추 씨'
$ws.Cells.Item(103, "B").Value = '옥수수'
$ws.Cells.Item(104, "B").Value = '옥수수 씨'
$ws.Cells.Item(105, "B").Value = '노란 파프리카'
$ws.Cells.Item(106, "B").Value = '노란 파프리카 씨'
$ws.Cells.Item(107, "B").Value = '초록 파프리카'
$ws.Cells.Item(108, "B").Value = '초록 파프리카 씨'
$ws.Cells.Item(109, "B").Value = '빨간 파프리카'
$ws.Cells.Item(110, "B").Value = '빨간 파프리카 씨'
$ws.Cells.Item(111, "B").Value = '용과'
$ws.Cells.Item(112, "B").Value = '용과 씨'
$ws.Cells.Item(113, "B").Value = '수박'
$ws.Cells.Item(114, "B").Value = '수박 씨'
$ws.Cells.Item(115, "B").Value = '오이'
$ws.Cells.Item(116, "B").Value = '오이 씨'
$ws.Cells.Item(117, "B").Value = '가지'
$ws.Cells.Item(118, "B").Value = '가지 씨'
$ws.Cells.Item(119, "B").Value = '파인애플'
$ws.Cells.Item(120, "B").Value = '파인애플 씨'
$ws.Cells.Item(121, "B").Value = '완두콩'
$ws.Cells.Item(122, "B").Value = '완두콩 씨'
$ws.Cells.Item(123, "B").Value = '파란 작물'
$ws.Cells.Item(124, "B").Value = '파란 작물 씨'
$ws.Cells.Item(125, "B").Value = '쌀'
$ws.Cells.Item(126, "B").Value = '쌀 씨'
$ws.Cells.Item(97, "L").Value = 'Tomato'
$ws.Cells.Item(99, "L").Value = 'Sunflower'
$ws.Cells.Item(101, "L").Value = 'Pepper'
$ws.Cells.Item(103, "L").Value = 'Corn'
$ws.Cells.Item(105, "L").Value = 'YelloPaprika'
$ws.Cells.Item(107, "L").Value = 'GreenPaprika'
$ws.Cells.Item(109, "L").Value = 'RedPaprika'
$ws.Cells.Item(111, "L").Value = 'Dragonfruit'
$ws.Cells.Item(113, "L").Value = 'Watermelon'
$ws.Cells.Item(115, "L").Value = 'Cucumber'
$ws.Cells.Item(117, "L").Value = 'Eggplant'
$ws.Cells.Item(119, "L").Value = 'Pineapple'
$ws.Cells.Item(121, "L").Value = 'Pea'
$ws.Cells.Item(123, "L").Value = 'Bluecrops'
$ws.Cells.Item(125, "L").Value = 'Rice'
$ws.Cells.Item(98, "L").Value = 'Tomato_Seed'
$ws.Cells.Item(100, "L").Value = 'Sunflower_Seed'
$ws.Cells.Item(102, "L").Value = 'Pepper_Seed'
$ws.Cells.Item(104, "L").Value = 'Corn_Seed'
$ws.Cells.Item(106, "L").Value = 'YelloPaprika_Seed'
$ws.Cells.Item(108, "L").Value = 'GreenPaprika_Seed'
$ws.Cells.Item(110, "L").Value = 'RedPaprika_Seed'
$ws.Cells.Item(112, "L").Value = 'Dragonfruit_Seed'
$ws.Cells.Item(114, "L").Value = 'Watermelon_Seed'
$ws.Cells.Item(116, "L").Value = 'Cucumber_Seed'
$ws.Cells.Item(118, "L").Value = 'Eggplant_Seed'
$ws.Cells.Item(120, "L").Value = 'Pineapple_Seed'
$ws.Cells.Item(122, "L").Value = 'Pea_Seed'
$ws.Cells.Item(124, "L").Value = 'Bluecrops_Seed'
$ws.Cells.Item(126, "L").Value = 'Rice_Seed'

# --- Phase 2: fill in the remaining cells for the 52 new rows (75-126) ---
$row = 75
$ws.Cells.Item($row, "A").Value = 74
$ws.Cells.Item($row, "C").Value = '딸기'
$ws.Cells.Item($row, "D").Value = 'Except'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 76
$ws.Cells.Item($row, "A").Value = 75
$ws.Cells.Item($row, "C").Value = '딸기 씨'
$ws.Cells.Item($row, "D").Value = 'Seed'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 77
$ws.Cells.Item($row, "A").Value = 76
$ws.Cells.Item($row, "C").Value = '대파'
$ws.Cells.Item($row, "D").Value = 'Except'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 78
$ws.Cells.Item($row, "A").Value = 77
$ws.Cells.Item($row, "C").Value = '대파 씨'
$ws.Cells.Item($row, "D").Value = 'Seed'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 79
$ws.Cells.Item($row, "A").Value = 78
$ws.Cells.Item($row, "C").Value = '감자'
$ws.Cells.Item($row, "D").Value = 'Except'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 80
$ws.Cells.Item($row, "A").Value = 79
$ws.Cells.Item($row, "C").Value = '감자 씨'
$ws.Cells.Item($row, "D").Value = 'Seed'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 81
$ws.Cells.Item($row, "A").Value = 80
$ws.Cells.Item($row, "C").Value = '양파'
$ws.Cells.Item($row, "D").Value = 'Except'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 82
$ws.Cells.Item($row, "A").Value = 81
$ws.Cells.Item($row, "C").Value = '양파 씨'
$ws.Cells.Item($row, "D").Value = 'Seed'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 83
$ws.Cells.Item($row, "A").Value = 82
$ws.Cells.Item($row, "C").Value = '당근'
$ws.Cells.Item($row, "D").Value = 'Except'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 84
$ws.Cells.Item($row, "A").Value = 83
$ws.Cells.Item($row, "C").Value = '당근 씨'
$ws.Cells.Item($row, "D").Value = 'Seed'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 85
$ws.Cells.Item($row, "A").Value = 84
$ws.Cells.Item($row, "C").Value = '블루베리'
$ws.Cells.Item($row, "D").Value = 'Except'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 86
$ws.Cells.Item($row, "A").Value = 85
$ws.Cells.Item($row, "C").Value = '블루베리 씨'
$ws.Cells.Item($row, "D").Value = 'Seed'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 87
$ws.Cells.Item($row, "A").Value = 86
$ws.Cells.Item($row, "C").Value = '무'
$ws.Cells.Item($row, "D").Value = 'Except'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 88
$ws.Cells.Item($row, "A").Value = 87
$ws.Cells.Item($row, "C").Value = '무 씨'
$ws.Cells.Item($row, "D").Value = 'Seed'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 89
$ws.Cells.Item($row, "A").Value = 88
$ws.Cells.Item($row, "C").Value = '양배추'
$ws.Cells.Item($row, "D").Value = 'Except'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 90
$ws.Cells.Item($row, "A").Value = 89
$ws.Cells.Item($row, "C").Value = '양배추 씨'
$ws.Cells.Item($row, "D").Value = 'Seed'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 91
$ws.Cells.Item($row, "A").Value = 90
$ws.Cells.Item($row, "C").Value = '콜리플라워'
$ws.Cells.Item($row, "D").Value = 'Except'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 92
$ws.Cells.Item($row, "A").Value = 91
$ws.Cells.Item($row, "C").Value = '콜리플라워 씨'
$ws.Cells.Item($row, "D").Value = 'Seed'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 93
$ws.Cells.Item($row, "A").Value = 92
$ws.Cells.Item($row, "C").Value = '밀'
$ws.Cells.Item($row, "D").Value = 'Except'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 94
$ws.Cells.Item($row, "A").Value = 93
$ws.Cells.Item($row, "C").Value = '밀 씨'
$ws.Cells.Item($row, "D").Value = 'Seed'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 95
$ws.Cells.Item($row, "A").Value = 94
$ws.Cells.Item($row, "C").Value = '브로콜리'
$ws.Cells.Item($row, "D").Value = 'Except'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 96
$ws.Cells.Item($row, "A").Value = 95
$ws.Cells.Item($row, "C").Value = '브로콜리 씨'
$ws.Cells.Item($row, "D").Value = 'Seed'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 97
$ws.Cells.Item($row, "A").Value = 96
$ws.Cells.Item($row, "C").Value = '토마토'
$ws.Cells.Item($row, "D").Value = 'Except'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 98
$ws.Cells.Item($row, "A").Value = 97
$ws.Cells.Item($row, "C").Value = '토마토 씨'
$ws.Cells.Item($row, "D").Value = 'Seed'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 99
$ws.Cells.Item($row, "A").Value = 98
$ws.Cells.Item($row, "C").Value = '해바라기'
$ws.Cells.Item($row, "D").Value = 'Except'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 100
$ws.Cells.Item($row, "A").Value = 99
$ws.Cells.Item($row, "C").Value = '해바라기 씨'
$ws.Cells.Item($row, "D").Value = 'Seed'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 101
$ws.Cells.Item($row, "A").Value = 100
$ws.Cells.Item($row, "C").Value = '고추'
$ws.Cells.Item($row, "D").Value = 'Except'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 102
$ws.Cells.Item($row, "A").Value = 101
$ws.Cells.Item($row, "C").Value = '고추 씨'
$ws.Cells.Item($row, "D").Value = 'Seed'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 103
$ws.Cells.Item($row, "A").Value = 102
$ws.Cells.Item($row, "C").Value = '옥수수'
$ws.Cells.Item($row, "D").Value = 'Except'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 104
$ws.Cells.Item($row, "A").Value = 103
$ws.Cells.Item($row, "C").Value = '옥수수 씨'
$ws.Cells.Item($row, "D").Value = 'Seed'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 105
$ws.Cells.Item($row, "A").Value = 104
$ws.Cells.Item($row, "C").Value = '노란 파프리카'
$ws.Cells.Item($row, "D").Value = 'Except'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 106
$ws.Cells.Item($row, "A").Value = 105
$ws.Cells.Item($row, "C").Value = '노란 파프리카 씨'
$ws.Cells.Item($row, "D").Value = 'Seed'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 107
$ws.Cells.Item($row, "A").Value = 106
$ws.Cells.Item($row, "C").Value = '초록 파프리카'
$ws.Cells.Item($row, "D").Value = 'Except'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 108
$ws.Cells.Item($row, "A").Value = 107
$ws.Cells.Item($row, "C").Value = '초록 파프리카 씨'
$ws.Cells.Item($row, "D").Value = 'Seed'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 109
$ws.Cells.Item($row, "A").Value = 108
$ws.Cells.Item($row, "C").Value = '빨간 파프리카'
$ws.Cells.Item($row, "D").Value = 'Except'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 110
$ws.Cells.Item($row, "A").Value = 109
$ws.Cells.Item($row, "C").Value = '빨간 파프리카 씨'
$ws.Cells.Item($row, "D").Value = 'Seed'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 111
$ws.Cells.Item($row, "A").Value = 110
$ws.Cells.Item($row, "C").Value = '용과'
$ws.Cells.Item($row, "D").Value = 'Except'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 112
$ws.Cells.Item($row, "A").Value = 111
$ws.Cells.Item($row, "C").Value = '용과 씨'
$ws.Cells.Item($row, "D").Value = 'Seed'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 113
$ws.Cells.Item($row, "A").Value = 112
$ws.Cells.Item($row, "C").Value = '수박'
$ws.Cells.Item($row, "D").Value = 'Except'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 114
$ws.Cells.Item($row, "A").Value = 113
$ws.Cells.Item($row, "C").Value = '수박 씨'
$ws.Cells.Item($row, "D").Value = 'Seed'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 115
$ws.Cells.Item($row, "A").Value = 114
$ws.Cells.Item($row, "C").Value = '오이'
$ws.Cells.Item($row, "D").Value = 'Except'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 116
$ws.Cells.Item($row, "A").Value = 115
$ws.Cells.Item($row, "C").Value = '오이 씨'
$ws.Cells.Item($row, "D").Value = 'Seed'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 117
$ws.Cells.Item($row, "A").Value = 116
$ws.Cells.Item($row, "C").Value = '가지'
$ws.Cells.Item($row, "D").Value = 'Except'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 118
$ws.Cells.Item($row, "A").Value = 117
$ws.Cells.Item($row, "C").Value = '가지 씨'
$ws.Cells.Item($row, "D").Value = 'Seed'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 119
$ws.Cells.Item($row, "A").Value = 118
$ws.Cells.Item($row, "C").Value = '파인애플'
$ws.Cells.Item($row, "D").Value = 'Except'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 120
$ws.Cells.Item($row, "A").Value = 119
$ws.Cells.Item($row, "C").Value = '파인애플 씨'
$ws.Cells.Item($row, "D").Value = 'Seed'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 121
$ws.Cells.Item($row, "A").Value = 120
$ws.Cells.Item($row, "C").Value = '완두콩'
$ws.Cells.Item($row, "D").Value = 'Except'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 122
$ws.Cells.Item($row, "A").Value = 121
$ws.Cells.Item($row, "C").Value = '완두콩 씨'
$ws.Cells.Item($row, "D").Value = 'Seed'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 123
$ws.Cells.Item($row, "A").Value = 122
$ws.Cells.Item($row, "C").Value = '파란 작물'
$ws.Cells.Item($row, "D").Value = 'Except'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 124
$ws.Cells.Item($row, "A").Value = 123
$ws.Cells.Item($row, "C").Value = '파란 작물 씨'
$ws.Cells.Item($row, "D").Value = 'Seed'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 125
$ws.Cells.Item($row, "A").Value = 124
$ws.Cells.Item($row, "C").Value = '쌀'
$ws.Cells.Item($row, "D").Value = 'Except'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64
$row = 126
$ws.Cells.Item($row, "A").Value = 125
$ws.Cells.Item($row, "C").Value = '쌀 씨'
$ws.Cells.Item($row, "D").Value = 'Seed'
$ws.Cells.Item($row, "E").Value = 1
$ws.Cells.Item($row, "F").Value = 64

# --- Phase 3: dimension / view bookkeeping ---
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 115
$ws.Range("T8").Select()
